# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the per-locale report sheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 12:58:45"
$wsZhCn.Range("H2").Value = "2016-03-21 12:59:07"

# de-de sheet: row 2 datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 12:58:49"
$wsDeDe.Range("H2").Value = "2016-03-21 12:59:13"
